# Update templates and docs to fit new data model.
# On the "CollectionEvents" sheet the "ageMin" / "ageMax" columns are
# replaced by a single "ageGroups" column (the "ageMax" column is removed,
# and the remaining "ageMin" header is renamed to "ageGroups"); the
# "subcohorts" column shifts one position to the left as a result.
# (The "SourceTables" sheet references the shared strings "unitOfObservation"
# and "numberOfRows" which keep their values but will end up pointing at
# different shared-string indices once "ageMin"/"ageMax" are removed and
# "ageGroups" is appended - that happens automatically when the workbook is
# saved.)

$wb = $excel.ActiveWorkbook

$ceSheet = $wb.Worksheets.Item("CollectionEvents")

# Remove the "ageMax" column (column D) entirely - "subcohorts" (old column E)
# shifts left into column D.
$null = $ceSheet.Range("D1").EntireColumn.Delete()

# Rename the remaining former "ageMin" header (now column C) to "ageGroups".
$ceSheet.Range("C1").Value = "ageGroups"

# Match the narrower "best fit" width the new column ends up with.
$ceSheet.Columns.Item(3).ColumnWidth = 8.285714285714286

# Update the sheet's selection to the (now empty) next column, as in the
# authored workbook.
$null = $ceSheet.Range("E1").Select()

# Restore the originally active sheet/selection so the workbook's active tab
# is unchanged.
$srcSheet = $wb.Worksheets.Item("SourceTables")
$null = $srcSheet.Activate()
$null = $srcSheet.Range("F7").Select()
